$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 546.7765096666667
$ws.Range("H2").Value = 1640.329529
$ws.Range("I2").Value = 0.6285526459909564
$ws.Range("J2").Value = 0.6285526459909564
$ws.Range("M2").Value = 281.0920463333333
$ws.Range("N2").Value = 843.2761389999999
$ws.Range("O2").Value = 0.8291026083535286
$ws.Range("P2").Value = 0.8291026083535286
$ws.Range("Q2").Value = 153694.527989201
$ws.Range("R2").Value = 1383250.751902808
$ws.Range("S2").Value = 0.5211346382786139
$ws.Range("T2").Value = 0.5211346382786139
$ws.Range("G3").Value = 546.7765096666667
$ws.Range("H3").Value = 1640.329529
$ws.Range("I3").Value = 0.6285526459909564
$ws.Range("J3").Value = 0.6285526459909564
$ws.Range("O3").Value = 0.001324719879221983
$ws.Range("P3").Value = 0.001324719879221983
$ws.Range("Q3").Value = 245.5693595745127
$ws.Range("R3").Value = 2210.124236170614
$ws.Range("S3").Value = 0.0008326561852817979
$ws.Range("T3").Value = 0.0008326561852817979
$ws.Range("G4").Value = 546.7765096666667
$ws.Range("H4").Value = 1640.329529
$ws.Range("I4").Value = 0.6285526459909564
$ws.Range("J4").Value = 0.6285526459909564
$ws.Range("M4").Value = 4.452417
$ws.Range("N4").Value = 13.357251
$ws.Range("O4").Value = 0.01313274635953239
$ws.Range("P4").Value = 0.01313274635953239
$ws.Range("Q4").Value = 2434.477026840531
$ws.Range("R4").Value = 21910.29324156478
$ws.Range("S4").Value = 0.008254622473412182
$ws.Range("T4").Value = 0.008254622473412182
$ws.Range("G5").Value = 546.7765096666667
$ws.Range("H5").Value = 1640.329529
$ws.Range("I5").Value = 0.6285526459909564
$ws.Range("J5").Value = 0.6285526459909564
$ws.Range("M5").Value = 53.03808999999999
$ws.Range("N5").Value = 159.11427
$ws.Range("O5").Value = 0.156439925407717
$ws.Range("P5").Value = 0.156439925407717
$ws.Range("Q5").Value = 28999.98172958653
$ws.Range("R5").Value = 260999.8355662788
$ws.Range("S5").Value = 0.09833072905364836
$ws.Range("T5").Value = 0.09833072905364837
$ws.Range("I6").Value = 0.1861770314550556
$ws.Range("J6").Value = 0.1861770314550556
$ws.Range("M6").Value = 281.0920463333333
$ws.Range("N6").Value = 843.2761389999999
$ws.Range("O6").Value = 0.8291026083535286
$ws.Range("P6").Value = 0.8291026083535286
$ws.Range("Q6").Value = 45524.25505552179
$ws.Range("R6").Value = 409718.2954996961
$ws.Range("S6").Value = 0.1543598623949035
$ws.Range("T6").Value = 0.1543598623949035
$ws.Range("I7").Value = 0.1861770314550556
$ws.Range("J7").Value = 0.1861770314550556
$ws.Range("O7").Value = 0.001324719879221983
$ws.Range("P7").Value = 0.001324719879221983
$ws.Range("S7").Value = 0.0002466324146230487
$ws.Range("T7").Value = 0.0002466324146230487
$ws.Range("I8").Value = 0.1861770314550556
$ws.Range("J8").Value = 0.1861770314550556
$ws.Range("M8").Value = 4.452417
$ws.Range("N8").Value = 13.357251
$ws.Range("O8").Value = 0.01313274635953239
$ws.Range("P8").Value = 0.01313274635953239
$ws.Range("Q8").Value = 721.091079472158
$ws.Range("R8").Value = 6489.819715249421
$ws.Range("S8").Value = 0.002445015732069929
$ws.Range("T8").Value = 0.002445015732069929
$ws.Range("I9").Value = 0.1861770314550556
$ws.Range("J9").Value = 0.1861770314550556
$ws.Range("M9").Value = 53.03808999999999
$ws.Range("N9").Value = 159.11427
$ws.Range("O9").Value = 0.156439925407717
$ws.Range("P9").Value = 0.156439925407717
$ws.Range("Q9").Value = 8589.782486959659
$ws.Range("R9").Value = 77308.04238263692
$ws.Range("S9").Value = 0.02912552091345908
$ws.Range("T9").Value = 0.02912552091345908
$ws.Range("G10").Value = 160.630483
$ws.Range("H10").Value = 481.891449
$ws.Range("I10").Value = 0.1846544489960017
$ws.Range("J10").Value = 0.1846544489960017
$ws.Range("M10").Value = 281.0920463333333
$ws.Range("N10").Value = 843.2761389999999
$ws.Range("O10").Value = 0.8291026083535286
$ws.Range("P10").Value = 0.8291026083535286
$ws.Range("Q10").Value = 45151.95116998171
$ws.Range("R10").Value = 406367.5605298353
$ws.Range("S10").Value = 0.1530974853066686
$ws.Range("T10").Value = 0.1530974853066686
$ws.Range("G11").Value = 160.630483
$ws.Range("H11").Value = 481.891449
$ws.Range("I11").Value = 0.1846544489960017
$ws.Range("J11").Value = 0.1846544489960017
$ws.Range("O11").Value = 0.001324719879221983
$ws.Range("P11").Value = 0.001324719879221983
$ws.Range("Q11").Value = 72.142683785926
$ws.Range("R11").Value = 649.284154073334
$ws.Range("S11").Value = 0.0002446154193717853
$ws.Range("T11").Value = 0.0002446154193717853
$ws.Range("G12").Value = 160.630483
$ws.Range("H12").Value = 481.891449
$ws.Range("I12").Value = 0.1846544489960017
$ws.Range("J12").Value = 0.1846544489960017
$ws.Range("M12").Value = 4.452417
$ws.Range("N12").Value = 13.357251
$ws.Range("O12").Value = 0.01313274635953239
$ws.Range("P12").Value = 0.01313274635953239
$ws.Range("Q12").Value = 715.193893227411
$ws.Range("R12").Value = 6436.745039046698
$ws.Range("S12").Value = 0.0024250200428237
$ws.Range("T12").Value = 0.0024250200428237
$ws.Range("G13").Value = 160.630483
$ws.Range("H13").Value = 481.891449
$ws.Range("I13").Value = 0.1846544489960017
$ws.Range("J13").Value = 0.1846544489960017
$ws.Range("M13").Value = 53.03808999999999
$ws.Range("N13").Value = 159.11427
$ws.Range("O13").Value = 0.156439925407717
$ws.Range("P13").Value = 0.156439925407717
$ws.Range("Q13").Value = 8519.534014097468
$ws.Range("R13").Value = 76675.80612687721
$ws.Range("S13").Value = 0.02888732822713758
$ws.Range("T13").Value = 0.02888732822713759
$ws.Range("G14").Value = 0.5357470000000001
$ws.Range("H14").Value = 1.607241
$ws.Range("I14").Value = 0.0006158735579862568
$ws.Range("J14").Value = 0.0006158735579862568
$ws.Range("M14").Value = 281.0920463333333
$ws.Range("N14").Value = 843.2761389999999
$ws.Range("O14").Value = 0.8291026083535286
$ws.Range("P14").Value = 0.8291026083535286
$ws.Range("Q14").Value = 150.5942205469443
$ws.Range("R14").Value = 1355.347984922499
$ws.Range("S14").Value = 0.0005106223733423736
$ws.Range("T14").Value = 0.0005106223733423736
$ws.Range("G15").Value = 0.5357470000000001
$ws.Range("H15").Value = 1.607241
$ws.Range("I15").Value = 0.0006158735579862568
$ws.Range("J15").Value = 0.0006158735579862568
$ws.Range("O15").Value = 0.001324719879221983
$ws.Range("P15").Value = 0.001324719879221983
$ws.Range("Q15").Value = 0.240615764134
$ws.Range("R15").Value = 2.165541877206
$ws.Range("S15").Value = 0.0000008158599453515673
$ws.Range("T15").Value = 0.0000008158599453515673
$ws.Range("G16").Value = 0.5357470000000001
$ws.Range("H16").Value = 1.607241
$ws.Range("I16").Value = 0.0006158735579862568
$ws.Range("J16").Value = 0.0006158735579862568
$ws.Range("M16").Value = 4.452417
$ws.Range("N16").Value = 13.357251
$ws.Range("O16").Value = 0.01313274635953239
$ws.Range("P16").Value = 0.01313274635953239
$ws.Range("Q16").Value = 2.385369050499
$ws.Range("R16").Value = 21.468321454491
$ws.Range("S16").Value = 0.000008088111226576272
$ws.Range("T16").Value = 0.000008088111226576273
$ws.Range("G17").Value = 0.5357470000000001
$ws.Range("H17").Value = 1.607241
$ws.Range("I17").Value = 0.0006158735579862568
$ws.Range("J17").Value = 0.0006158735579862568
$ws.Range("M17").Value = 53.03808999999999
$ws.Range("N17").Value = 159.11427
$ws.Range("O17").Value = 0.156439925407717
$ws.Range("P17").Value = 0.156439925407717
$ws.Range("Q17").Value = 28.41499760323
$ws.Range("R17").Value = 255.73497842907
$ws.Range("S17").Value = 0.00009634721347195527
$ws.Range("T17").Value = 0.00009634721347195529

Write-Host "Updated 174 cells"
